# Master data update: add UIN deactivation/reactivation + registration
# acknowledgement template rows, and drop the now-unused "Sheet1" helper
# sheet (and its Table1) that used to hold the filtered/hidden view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-template_type")

# Remove the helper "Sheet1" (it carried Table1 / the autofiltered copy of
# the data) - it is no longer needed.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Delete()

$newRows = @(
    @{Row=92; A="RPR_UIN_CARD_TEMPLATE"; B="UIN card template"; C="eng"},
    @{Row=93; A="RPR_UIN_CARD_TEMPLATE"; B="قالب بطاقة UIN"; C="ara"},
    @{Row=94; A="RPR_UIN_CARD_TEMPLATE"; B="Modèle de carte UIN"; C="fra"},
    @{Row=95; A="RPR_UIN_DEAC_SMS"; B="Template for UIN Deactivation SMS"; C="eng"},
    @{Row=96; A="RPR_UIN_DEAC_SMS"; B="قالب لتعطيل UIN SMS"; C="ara"},
    @{Row=97; A="RPR_UIN_DEAC_SMS"; B="Modèle pour SMS de désactivation UIN"; C="fra"},
    @{Row=98; A="RPR_UIN_DEAC_EMAIL"; B="Template for UIN Deactivation Email"; C="eng"},
    @{Row=99; A="RPR_UIN_DEAC_EMAIL"; B="قالب لإلغاء تنشيط البريد"; C="ara"},
    @{Row=100; A="RPR_UIN_DEAC_EMAIL"; B="Modèle pour Email de désactivation UIN"; C="fra"},
    @{Row=101; A="RPR_UIN_REAC_SMS"; B="Template for UIN Reactivate SMS"; C="eng"},
    @{Row=102; A="RPR_UIN_REAC_SMS"; B="قالب لـ UIN تنشيط SMS"; C="ara"},
    @{Row=103; A="RPR_UIN_REAC_SMS"; B="Modèle pour UIN Réactiver SMS"; C="fra"},
    @{Row=104; A="RPR_UIN_REAC_EMAIL"; B="Template for UIN Reactivate Email"; C="eng"},
    @{Row=105; A="RPR_UIN_REAC_EMAIL"; B="قالب لـ UIN تنشيط البريد"; C="ara"},
    @{Row=106; A="RPR_UIN_REAC_EMAIL"; B="Modèle pour UIN Réactiver Email"; C="fra"},
    @{Row=107; A="reg-sms-notification"; B="Registration Acknowledgement Template"; C="eng"},
    @{Row=108; A="reg-sms-notification"; B="نموذج شكر التسجيل"; C="ara"},
    @{Row=109; A="reg-sms-notification"; B="accusé de réception"; C="fra"},
    @{Row=110; A="reg-email-notification"; B="Registration Acknowledgement Template"; C="eng"},
    @{Row=111; A="reg-email-notification"; B="نموذج شكر التسجيل"; C="ara"},
    @{Row=112; A="reg-email-notification"; B="accusé de réception"; C="fra"},
    @{Row=113; A="reg-ack-template-part1"; B="Registration Acknowledgement Template - Part 1"; C="eng"},
    @{Row=114; A="reg-ack-template-part2"; B="نموذج شكر التسجيل"; C="ara"},
    @{Row=115; A="reg-ack-template-part3"; B="accusé de réception"; C="fra"},
    @{Row=116; A="reg-ack-template-part2"; B="Registration Acknowledgement Template - Part 2"; C="eng"},
    @{Row=117; A="reg-ack-template-part3"; B="نموذج شكر التسجيل"; C="ara"},
    @{Row=118; A="reg-ack-template-part4"; B="accusé de réception"; C="fra"},
    @{Row=119; A="reg-ack-template-part3"; B="Registration Acknowledgement Template - Part 3"; C="eng"},
    @{Row=120; A="reg-ack-template-part4"; B="نموذج شكر التسجيل"; C="ara"},
    @{Row=121; A="reg-ack-template-part5"; B="accusé de réception"; C="fra"},
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $true
    $ws.Cells.Item($r.Row, 5).Value = "superadmin"
    $ws.Cells.Item($r.Row, 6).Value = "now()"
}

# Move the selection below the new data, matching where the user's cursor
# ended up after pasting the new rows (selecting the rest of the sheet).
$ws.Range("A122:XFD1048576").Select()

Write-Host "Applied master-template_type updates"
